$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 25.00614441910968
$ws.Cells.Item(2, 3).Value = 12.30731053561339
$ws.Cells.Item(2, 4).Value = 4.060659819108812
$ws.Cells.Item(2, 5).Value = 9.760591166727076
$ws.Cells.Item(2, 6).Value = 53.83006190449456
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 9).Value = 37.56898058156668
$ws.Cells.Item(2, 10).Value = 9.65349555723383
$ws.Cells.Item(2, 12).Value = 12.678561840459
$ws.Cells.Item(3, 2).Value = 24.74140692267946
$ws.Cells.Item(3, 3).Value = 11.95175431647435
$ws.Cells.Item(3, 4).Value = 4.034675934778834
$ws.Cells.Item(3, 5).Value = 9.769613835736122
$ws.Cells.Item(3, 6).Value = 53.70344327964852
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 9).Value = 37.55785225907313
$ws.Cells.Item(3, 10).Value = 9.670375217775941
$ws.Cells.Item(3, 12).Value = 12.68597716097472
$ws.Cells.Item(4, 2).Value = 24.58469654994455
$ws.Cells.Item(4, 3).Value = 11.73228476847309
$ws.Cells.Item(4, 4).Value = 4.018360675747212
$ws.Cells.Item(4, 5).Value = 9.775533184945656
$ws.Cells.Item(4, 6).Value = 53.6380974084358
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 9).Value = 37.55859009833727
$ws.Cells.Item(4, 10).Value = 9.681314013636113
$ws.Cells.Item(4, 12).Value = 12.69273076665188
$ws.Cells.Item(5, 2).Value = 24.5223769275508
$ws.Cells.Item(5, 3).Value = 11.64271556975238
$ws.Cells.Item(5, 4).Value = 4.011621337034558
$ws.Cells.Item(5, 5).Value = 9.778040985024914
$ws.Cells.Item(5, 6).Value = 53.61459065031846
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 9).Value = 37.5607887618199
$ws.Cells.Item(5, 10).Value = 9.685916554616993
$ws.Cells.Item(5, 12).Value = 12.69603650654151
$ws.Cells.Item(6, 2).Value = 24.51212394373673
$ws.Cells.Item(6, 3).Value = 11.62783935030638
$ws.Cells.Item(6, 4).Value = 4.010496795684405
$ws.Cells.Item(6, 5).Value = 9.778463184749462
$ws.Cells.Item(6, 6).Value = 53.61087601526252
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 9).Value = 37.56126829473202
$ws.Cells.Item(6, 10).Value = 9.686689566000526
$ws.Cells.Item(6, 12).Value = 12.69661886155735
$ws.Cells.Item(7, 2).Value = 24.58384975219032
$ws.Cells.Item(7, 3).Value = 11.73107713427301
$ws.Cells.Item(7, 4).Value = 4.018270153979174
$ws.Cells.Item(7, 5).Value = 9.775566618563897
$ws.Cells.Item(7, 6).Value = 53.6377677430326
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 9).Value = 37.55861207394265
$ws.Cells.Item(7, 10).Value = 9.681375497920248
$ws.Cells.Item(7, 12).Value = 12.69277310736525
$ws.Cells.Item(8, 2).Value = 24.91369552549676
$ws.Cells.Item(8, 3).Value = 12.18504667118562
$ws.Cells.Item(8, 4).Value = 4.051775543143769
$ws.Cells.Item(8, 5).Value = 9.763623568833102
$ws.Cells.Item(8, 6).Value = 53.78383412554785
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 9).Value = 37.56356909335558
$ws.Cells.Item(8, 10).Value = 9.659196644569688
$ws.Cells.Item(8, 12).Value = 12.68066210169591
$ws.Cells.Item(9, 2).Value = 25.60351802926547
$ws.Cells.Item(9, 3).Value = 13.05962207465613
$ws.Cells.Item(9, 4).Value = 4.114613666486785
$ws.Cells.Item(9, 5).Value = 9.743204044116675
$ws.Cells.Item(9, 6).Value = 54.16831647701361
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 9).Value = 37.63357610836246
$ws.Cells.Item(9, 10).Value = 9.62024497442753
$ws.Cells.Item(9, 12).Value = 12.67435775374522
$ws.Cells.Item(10, 2).Value = 26.13191586584594
$ws.Cells.Item(10, 3).Value = 13.68462226798599
$ws.Cells.Item(10, 4).Value = 4.159034210943485
$ws.Cells.Item(10, 5).Value = 9.730018175672992
$ws.Cells.Item(10, 6).Value = 54.50993715559531
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 9).Value = 37.72200093700716
$ws.Cells.Item(10, 10).Value = 9.594369800315228
$ws.Cells.Item(10, 12).Value = 12.6803336106378
$ws.Cells.Item(11, 2).Value = 26.37597358210042
$ws.Cells.Item(11, 3).Value = 13.963652488355
$ws.Cells.Item(11, 4).Value = 4.178862930283786
$ws.Cells.Item(11, 5).Value = 9.724411289378757
$ws.Cells.Item(11, 6).Value = 54.67799241578483
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 9).Value = 37.77028858014415
$ws.Cells.Item(11, 10).Value = 9.583188576958577
$ws.Cells.Item(11, 12).Value = 12.68534633564183
$ws.Cells.Item(12, 2).Value = 26.46883680303258
$ws.Cells.Item(12, 3).Value = 14.06844431505563
$ws.Cells.Item(12, 4).Value = 4.186317023389177
$ws.Cells.Item(12, 5).Value = 9.722344190777147
$ws.Cells.Item(12, 6).Value = 54.74342764874384
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 9).Value = 37.78973349258629
$ws.Cells.Item(12, 10).Value = 9.579038897800963
$ws.Cells.Item(12, 12).Value = 12.68757331588333
$ws.Cells.Item(13, 2).Value = 26.44881867511234
$ws.Cells.Item(13, 3).Value = 14.04591590541677
$ws.Cells.Item(13, 4).Value = 4.184714083251882
$ws.Cells.Item(13, 5).Value = 9.722786884873678
$ws.Cells.Item(13, 6).Value = 54.72925550072559
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 9).Value = 37.78549412880498
$ws.Cells.Item(13, 10).Value = 9.579928857197201
$ws.Cells.Item(13, 12).Value = 12.68707909384386
$ws.Cells.Item(14, 2).Value = 26.38360509226414
$ws.Cells.Item(14, 3).Value = 13.97229182593235
$ws.Cells.Item(14, 4).Value = 4.179477278696127
$ws.Cells.Item(14, 5).Value = 9.724240104468867
$ws.Cells.Item(14, 6).Value = 54.68333995098673
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 9).Value = 37.77186509211851
$ws.Cells.Item(14, 10).Value = 9.582845490788619
$ws.Cells.Item(14, 12).Value = 12.6855229688011
$ws.Cells.Item(15, 2).Value = 26.3437151395019
$ws.Cells.Item(15, 3).Value = 13.92707842933739
$ws.Cells.Item(15, 4).Value = 4.176262466475327
$ws.Cells.Item(15, 5).Value = 9.725137545022843
$ws.Cells.Item(15, 6).Value = 54.65544856359661
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 9).Value = 37.76366788634284
$ws.Cells.Item(15, 10).Value = 9.584642994626581
$ws.Cells.Item(15, 12).Value = 12.68461257580505
$ws.Cells.Item(16, 2).Value = 26.11603292349054
$ws.Cells.Item(16, 3).Value = 13.66627091675533
$ws.Cells.Item(16, 4).Value = 4.157730726296388
$ws.Cells.Item(16, 5).Value = 9.73039246002195
$ws.Cells.Item(16, 6).Value = 54.49920702922173
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 9).Value = 37.71900749689054
$ws.Cells.Item(16, 10).Value = 9.595112350987584
$ws.Cells.Item(16, 12).Value = 12.68005209031179
$ws.Cells.Item(17, 2).Value = 25.97724005637678
$ws.Cells.Item(17, 3).Value = 13.50483802759895
$ws.Cells.Item(17, 4).Value = 4.146265047277029
$ws.Cells.Item(17, 5).Value = 9.733716305990226
$ws.Cells.Item(17, 6).Value = 54.40658272081218
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 9).Value = 37.69367523640297
$ws.Cells.Item(17, 10).Value = 9.601685692304565
$ws.Cells.Item(17, 12).Value = 12.67784124289607
$ws.Cells.Item(18, 2).Value = 25.89776225969422
$ws.Cells.Item(18, 3).Value = 13.41149455387268
$ws.Cells.Item(18, 4).Value = 4.139634686927287
$ws.Cells.Item(18, 5).Value = 9.735664948171181
$ws.Cells.Item(18, 6).Value = 54.35449992927515
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 9).Value = 37.67986376678986
$ws.Cells.Item(18, 10).Value = 9.605522014766111
$ws.Cells.Item(18, 12).Value = 12.67678565019505
$ws.Cells.Item(19, 2).Value = 25.87091556736424
$ws.Cells.Item(19, 3).Value = 13.37980918282607
$ws.Cells.Item(19, 4).Value = 4.137383651496298
$ws.Cells.Item(19, 5).Value = 9.736331060165448
$ws.Cells.Item(19, 6).Value = 54.33707100681742
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 9).Value = 37.67531777583559
$ws.Cells.Item(19, 10).Value = 9.606830472762427
$ws.Cells.Item(19, 12).Value = 12.67646537958952
$ws.Cells.Item(20, 2).Value = 25.99197897019615
$ws.Cells.Item(20, 3).Value = 13.52207454336581
$ws.Cells.Item(20, 4).Value = 4.147489274761808
$ws.Cells.Item(20, 5).Value = 9.73335866398982
$ws.Cells.Item(20, 6).Value = 54.41631949744021
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 9).Value = 37.69629334380553
$ws.Cells.Item(20, 10).Value = 9.600980206440767
$ws.Cells.Item(20, 12).Value = 12.67805424099444
$ws.Cells.Item(21, 2).Value = 26.4027485416209
$ws.Cells.Item(21, 3).Value = 13.99394146899034
$ws.Cells.Item(21, 4).Value = 4.181016938938058
$ws.Cells.Item(21, 5).Value = 9.723811737318872
$ws.Cells.Item(21, 6).Value = 54.69677790434425
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 9).Value = 37.77583680956083
$ws.Cells.Item(21, 10).Value = 9.581986516905914
$ws.Cells.Item(21, 12).Value = 12.68597112773582
$ws.Cells.Item(22, 2).Value = 26.67375777348125
$ws.Cells.Item(22, 3).Value = 14.29721364333665
$ws.Cells.Item(22, 4).Value = 4.202610782909631
$ws.Cells.Item(22, 5).Value = 9.717899225924764
$ws.Cells.Item(22, 6).Value = 54.89053159704206
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 9).Value = 37.83458029917205
$ws.Cells.Item(22, 10).Value = 9.570064866450572
$ws.Cells.Item(22, 12).Value = 12.69306093623583
$ws.Cells.Item(23, 2).Value = 26.52891030994212
$ws.Cells.Item(23, 3).Value = 14.13585399160813
$ws.Cells.Item(23, 4).Value = 4.191114945756553
$ws.Cells.Item(23, 5).Value = 9.721024986340513
$ws.Cells.Item(23, 6).Value = 54.78617297665297
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 9).Value = 37.80260979872456
$ws.Cells.Item(23, 10).Value = 9.576382793757872
$ws.Cells.Item(23, 12).Value = 12.68910211313447
$ws.Cells.Item(24, 2).Value = 25.98531451621566
$ws.Cells.Item(24, 3).Value = 13.51428356797174
$ws.Cells.Item(24, 4).Value = 4.146935921696628
$ws.Cells.Item(24, 5).Value = 9.733520236438824
$ws.Cells.Item(24, 6).Value = 54.41191386048805
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 9).Value = 37.69510735395384
$ws.Cells.Item(24, 10).Value = 9.601298978341484
$ws.Cells.Item(24, 12).Value = 12.67795727338163
$ws.Cells.Item(25, 2).Value = 25.41277102371332
$ws.Cells.Item(25, 3).Value = 12.82556824623507
$ws.Cells.Item(25, 4).Value = 4.097920975809901
$ws.Cells.Item(25, 5).Value = 9.748408170706938
$ws.Cells.Item(25, 6).Value = 54.05385976446679
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 9).Value = 37.60815429206506
$ws.Cells.Item(25, 10).Value = 9.630298954151094
$ws.Cells.Item(25, 12).Value = 12.67419775927931
